$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomePage")

# Fill in new row 5 data (Element Name, Locator Strategy, Web Locator)
# Order matters for shared string table insertion order: search_query_top must be
# added before txt_search_bar to match target index order.
$ws.Range("C5").Value = "search_query_top"
$ws.Range("A5").Value = "txt_search_bar"
$ws.Range("B5").Value = "ID"

# Update the selection on the active sheet to G13
$ws.Activate()
$ws.Range("G13").Select()
